$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "46.894.93"
$ws.Range("E2").Value = "  +3.02%  "
# Row 3
Set-TextValue $ws.Range("D3") "2.628.88"
$ws.Range("E3").Value = "  +7.73%  "
# Row 4
$ws.Range("E4").Value = "  -0.32%  "
# Row 5
Set-TextValue $ws.Range("D5") "309.06"
$ws.Range("E5").Value = "  +5.37%  "
# Row 6
Set-TextValue $ws.Range("D6") "101.73"
$ws.Range("E6").Value = "  +7.68%  "
# Row 7
Set-TextValue $ws.Range("D7") "0.605"
$ws.Range("E7").Value = "  +6.45%  "
# Row 8
$ws.Range("E8").Value = "  +0.03%  "
# Row 9
Set-TextValue $ws.Range("D9") "0.584"
$ws.Range("E9").Value = "  +14.01%  "
# Row 10
Set-TextValue $ws.Range("D10") "39.95"
$ws.Range("E10").Value = "  +15.02%  "
# Row 11
Set-TextValue $ws.Range("D11") "0.0854"
$ws.Range("E11").Value = "  +9.57%  "
# Row 12
Set-TextValue $ws.Range("D12") "54.61"
$ws.Range("E12").Value = "  +1.92%  "
# Row 13
Set-TextValue $ws.Range("D13") "8.27"
$ws.Range("E13").Value = "  +14.09%  "
# Row 14
Set-TextValue $ws.Range("D14") "3.026.24"
$ws.Range("E14").Value = "  +7.00%  "
# Row 15
$ws.Range("E15").Value = "  +2.24%  "
# Row 16
Set-TextValue $ws.Range("D16") "2.629.83"
$ws.Range("E16").Value = "  +7.62%  "
# Row 17
Set-TextValue $ws.Range("D17") "0.936"
$ws.Range("E17").Value = "  +10.61%  "
# Row 18
Set-TextValue $ws.Range("D18") "15.08"
$ws.Range("E18").Value = "  +6.82%  "
# Row 19
Set-TextValue $ws.Range("D19") "46.998.71"
$ws.Range("E19").Value = "  +3.13%  "
# Row 20
Set-TextValue $ws.Range("D20") "0.0000102"
$ws.Range("E20").Value = "  +9.10%  "
# Row 21
Set-TextValue $ws.Range("D21") "13.09"
$ws.Range("E21").Value = "  +4.33%  "
# Row 22
Set-TextValue $ws.Range("D22") "6.79"
$ws.Range("E22").Value = "  +8.66%  "
# Row 23
Set-TextValue $ws.Range("D23") "277.35"
$ws.Range("E23").Value = "  +13.76%  "
# Row 24
Set-TextValue $ws.Range("D24") "72.01"
$ws.Range("E24").Value = "  +7.09%  "
# Row 25
Set-TextValue $ws.Range("D25") "3.06"
$ws.Range("E25").Value = "  +10.24%  "
# Row 26
Set-TextValue $ws.Range("D26") "2.19"
$ws.Range("E26").Value = "  +13.42%  "
# Row 27
$ws.Range("E27").Value = "  +37.32%  "
# Row 28
$ws.Range("E28").Value = "  +0.12%  "
# Row 29
$ws.Range("E29").Value = "  -0.19%  "
# Row 30
Set-TextValue $ws.Range("D30") "10.72"
$ws.Range("E30").Value = "  +10.07%  "
# Row 31
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D31") "2.32"
$ws.Range("E31").Value = "  +4.59%  "
# Row 32
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D32") "39.70"
$ws.Range("E32").Value = "  +2.20%  "
# Row 33
Set-TextValue $ws.Range("D33") "6.43"
$ws.Range("E33").Value = "  +17.00%  "
# Row 34
Set-TextValue $ws.Range("D34") "3.67"
$ws.Range("E34").Value = "  -2.65%  "
# Row 35
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D35") "2.89"
$ws.Range("E35").Value = "  +5.01%  "
# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D36") "2.27"
$ws.Range("E36").Value = "  +11.94%  "
# Row 37
Set-TextValue $ws.Range("D37") "0.0850"
$ws.Range("E37").Value = "  +11.12%  "
# Row 38
Set-TextValue $ws.Range("D38") "152.75"
$ws.Range("E38").Value = "  +4.51%  "
# Row 39
$ws.Range("E39").Value = "  +8.76%  "
# Row 40
$ws.Range("E40").Value = "  +7.30%  "
# Row 41
Set-TextValue $ws.Range("D41") "23.44"
$ws.Range("E41").Value = "  +44.63%  "
# Row 42
Set-TextValue $ws.Range("D42") "16.13"
$ws.Range("E42").Value = "  +8.46%  "
# Row 43
Set-TextValue $ws.Range("D43") "3.71"
$ws.Range("E43").Value = "  +15.60%  "
# Row 44
$ws.Range("E44").Value = "  +12.92%  "
# Row 45
Set-TextValue $ws.Range("D45") "4.14"
$ws.Range("E45").Value = "  +6.19%  "
# Row 46
Set-TextValue $ws.Range("D46") "2.133.27"
$ws.Range("E46").Value = "  +6.52%  "
# Row 47
Set-TextValue $ws.Range("D47") "0.998"
$ws.Range("E47").Value = "  -0.08%  "
# Row 48
Set-TextValue $ws.Range("D48") "94.13"
$ws.Range("E48").Value = "  +3.07%  "
# Row 49
Set-TextValue $ws.Range("D49") "9.61"
$ws.Range("E49").Value = "  +12.81%  "
# Row 50
$ws.Range("E50").Value = "  +3.49%  "
# Row 51
Set-TextValue $ws.Range("D51") "110.25"
$ws.Range("E51").Value = "  +7.89%  "
